# Updates cryptos list data cells per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.361.76'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '1.592.74'
$ws.Range('E3').Value = '  +0.46%  '
$ws.Range('E4').Value = '  -0.37%  '
$ws.Range('D5').Value = "'211.48"
$ws.Range('E5').Value = '  +0.84%  '
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('E7').Value = '  -0.34%  '
$ws.Range('E8').Value = '  +0.55%  '
$ws.Range('D10').Value = "'19.51"
$ws.Range('E10').Value = '  -0.48%  '
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('E13').Value = '  +1.27%  '
$ws.Range('D14').Value = '1.577.32'
$ws.Range('E14').Value = '  -1.67%  '
$ws.Range('D15').Value = "'0.524"
$ws.Range('D16').Value = "'64.72"
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').Value = '26.356.63'
$ws.Range('E17').Value = '  +0.29%  '
$ws.Range('D18').Value = '0.0₃0732'
$ws.Range('E18').Value = '  -1.07%  '
$ws.Range('D19').Value = "'7.52"
$ws.Range('E19').Value = '  +4.33%  '
$ws.Range('D20').Value = "'211.78"
$ws.Range('E20').Value = '  +2.34%  '
$ws.Range('E21').Value = '  -0.37%  '
$ws.Range('D23').Value = "'9.02"
$ws.Range('E23').Value = '  +2.13%  '
$ws.Range('E24').Value = '  -2.58%  '
$ws.Range('D25').Value = "'143.94"
$ws.Range('E25').Value = '  -0.35%  '
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('E27').Value = '  +0.66%  '
$ws.Range('E28').Value = '  -0.74%  '
$ws.Range('D29').Value = "'15.25"
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('E30').Value = '  -0.04%  '
$ws.Range('E31').Value = '  +0.91%  '
$ws.Range('E32').Value = '  -0.43%  '
$ws.Range('E33').Value = '  +1.45%  '
$ws.Range('D34').Value = '1.333.64'
$ws.Range('E34').Value = '  +3.88%  '
$ws.Range('E35').Value = '  -1.36%  '
$ws.Range('D36').Value = "'0.602"
$ws.Range('E36').Value = '  -1.44%  '
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('E38').Value = '  -0.38%  '
$ws.Range('D39').Value = "'0.819"
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').Value = "'5.76"
$ws.Range('E40').Value = '  +5.27%  '
$ws.Range('E41').Value = '  -0.32%  '
$ws.Range('D42').Value = "'1.00"
$ws.Range('E42').Value = '  -23.01%  '
$ws.Range('E43').Value = '  +0.19%  '
$ws.Range('E44').Value = '  -0.56%  '
$ws.Range('D45').Value = '1.729.64'
$ws.Range('E45').Value = '  +0.47%  '
$ws.Range('D46').Value = "'61.93"
$ws.Range('E46').Value = '  -0.80%  '
$ws.Range('D47').Value = "'88.20"
$ws.Range('E47').Value = '  -0.89%  '
$ws.Range('D48').Value = "'1.50"
$ws.Range('E48').Value = '  -3.85%  '
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = "'0.0982"
$ws.Range('E49').Value = '  -3.97%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = "'0.0504"
$ws.Range('E50').Value = '  -0.97%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = "'1.00"
$ws.Range('E51').Value = '  -0.39%  '
